# Apply the documented content updates to the generated FHIR StructureDefinition
# workbook (docs/StructureDefinition-observation-goal-reference.xlsx):
#   1. Metadata!B8  - "Date" bumped to the new publish timestamp
#   2. Metadata!B11 - "Description" trimmed (drop the trailing sentence)
#   3. Elements!K6  - Extension.value[x] Type(s) reference renamed from
#                      onc-patient-goal -> onc-nursing-goal
#   4. Elements column K widened slightly to keep "best fit" in sync with
#      the new (slightly wider-rendering) reference text

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2026-01-01T13:37:23+00:00"
$wsMeta.Range("B11").Value = "Extension to link goal evaluation observations to the patient goals being evaluated."

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("K6").Value = "Reference(https://clinyqai.github.io/open-nursing-core-ig/StructureDefinition/onc-nursing-goal)`n"

# Writing the (wrap-text) cell above makes the host recompute an explicit
# row height; put row 6 back to its original auto/default height so no
# stray ht/customHeight attributes are introduced.
$wsElem.Rows.Item(6).AutoFit()

# Nudge the column K "best fit" width to reflect the updated text metrics.
$wsElem.Columns.Item(11).ColumnWidth = 77.33333333333334
